$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.419591
$ws.Range("H2").Value = 4.258773
$ws.Range("I2").Value = 0.001848767113890483
$ws.Range("J2").Value = 0.001848767113890483
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 109.43559351829
$ws.Range("R2").Value = 984.92034166461
$ws.Range("S2").Value = 0.0004444097214207762
$ws.Range("T2").Value = 0.0004444097214207762
$ws.Range("G3").Value = 1.419591
$ws.Range("H3").Value = 4.258773
$ws.Range("I3").Value = 0.001848767113890483
$ws.Range("J3").Value = 0.001848767113890483
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 144.202106778064
$ws.Range("R3").Value = 1297.818961002576
$ws.Range("S3").Value = 0.0005855939191376335
$ws.Range("T3").Value = 0.0005855939191376334
$ws.Range("G4").Value = 1.419591
$ws.Range("H4").Value = 4.258773
$ws.Range("I4").Value = 0.001848767113890483
$ws.Range("J4").Value = 0.001848767113890483
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 201.619951896496
$ws.Range("R4").Value = 1814.579567068464
$ws.Range("S4").Value = 0.0008187634733320736
$ws.Range("T4").Value = 0.0008187634733320736
$ws.Range("I5").Value = 0.9578582377148513
$ws.Range("J5").Value = 0.9578582377148513
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 56699.29109141303
$ws.Range("R5").Value = 510293.6198227173
$ws.Range("S5").Value = 0.2302515602885551
$ws.Range("T5").Value = 0.2302515602885551
$ws.Range("I6").Value = 0.9578582377148513
$ws.Range("J6").Value = 0.9578582377148513
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.3034000092209203
$ws.Range("T6").Value = 0.3034000092209203
$ws.Range("I7").Value = 0.9578582377148513
$ws.Range("J7").Value = 0.9578582377148513
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.4242066682053759
$ws.Range("T7").Value = 0.4242066682053759
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.04029299517125823
$ws.Range("J8").Value = 0.04029299517125823
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2385.096428893663
$ws.Range("R8").Value = 21465.86786004297
$ws.Range("S8").Value = 0.009685697362706493
$ws.Range("T8").Value = 0.009685697362706493
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.04029299517125823
$ws.Range("J9").Value = 0.04029299517125823
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 3142.815960127462
$ws.Range("R9").Value = 28285.34364114716
$ws.Range("S9").Value = 0.01276273943800181
$ws.Range("T9").Value = 0.01276273943800181
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.04029299517125823
$ws.Range("J10").Value = 0.04029299517125823
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.01784455837054993
$ws.Range("T10").Value = 0.01784455837054993
